# Update "JSON Schema Restriction" sheet: insert a new "array" type row
# (B6:E6) between the existing "bigint[]" row (row 5) and "int" row
# (old row 6, now row 7), shifting every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6; everything from old row 6 onward
# shifts down to row 7 onward.
$ws.Rows.Item(6).Insert()

# Pick up the border/font/alignment formatting used by the other data rows:
#  - B6 (TYPE column, no wrap) + C6 (RESTRICTION column, wraps) come from
#    the row directly below (old "int" row, now row 7).
$ws.Range("B7:C7").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)  # xlPasteFormats

#  - D6/E6 (ACCEPTABLE VALUE / LENGTH columns, wraps) come from the
#    "bigint[]" row above (row 5), which is the nearest row that also has
#    those two columns populated/formatted.
$ws.Range("D5:E5").Copy()
$ws.Range("D6:E6").PasteSpecial(-4122)  # xlPasteFormats

# New cell text. Fill D6 before C6 so the shared-strings table gets the
# two new entries in the same order as the reference workbook (the
# bracket-array example string first, the JSON-schema snippet second).
$acceptableValues = "[]`n[123]`n[123,456,789]`n[ 123 , 456 , 789 ]"
$restriction = "`"type`": `"array`",`n`"items`": {`n    `"type`": `"number`"`n    }"

$ws.Cells.Item(6, 4).Value = $acceptableValues
$ws.Cells.Item(6, 3).Value = $restriction

# B6 (TYPE) and E6 (LENGTH) stay empty for this row.
$ws.Cells.Item(6, 2).Value = ""
$ws.Cells.Item(6, 5).Value = ""

# Match the taller row height used for other multi-line restriction rows.
$ws.Rows.Item(6).RowHeight = 51

# Move the active selection the way the authored workbook left it.
$ws.Range("D7").Select()
